# Moved statistical datasets and results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column A, shifting existing data (A:D) to (B:E)
$ws.Columns.Item(1).Insert()

# Insert a new row before row 1, shifting existing data (now rows 1:23) down to (2:24)
$ws.Rows.Item(1).Insert()

# Header row for the numeric columns (now B..E)
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# Row labels in the new column A, describing each pairwise comparison
$labels = @(
  "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)",
  "CyclomaticComplexity(CC) & NbOperators",
  "CyclomaticComplexity(CC) & EffortToImplement",
  "MaintainabilityIndex & MaintainabilityIndex",
  "NbUniqueOperands & NbUniqueOperands",
  "NbOperands & NbOperands",
  "NbOperands & EffortToImplement",
  "NbUniqueOperators & NbUniqueOperators",
  "NbOperators & CyclomaticComplexity(CC)",
  "NbOperators & NbOperators",
  "NbOperators & EffortToImplement",
  "ProgramLength & ProgramLength",
  "ProgramLength & EffortToImplement",
  "VocabularySize & VocabularySize",
  "ProgramVolume & ProgramVolume",
  "DifficultyLevel & DifficultyLevel",
  "ProgramLevel & ProgramLevel",
  "EffortToImplement & CyclomaticComplexity(CC)",
  "EffortToImplement & NbOperands",
  "EffortToImplement & NbOperators",
  "EffortToImplement & ProgramLength",
  "EffortToImplement & EffortToImplement",
  "TimeToImplement & TimeToImplement"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# Column width for the new label column A (columns B:E keep the widths
# that shifted over automatically from the original A:D on Insert()).
# 53.666667 is the closest input that snaps (via Excel's pixel grid) to
# the target stored width of ~54.552101.
$ws.Columns.Item(1).ColumnWidth = 53.666667
